$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# CH4 is now complete, so the "next up" marker moves from CH4 to CH5
$ws.Range("B3").Value = "Ch 5 - START"
